# Apply the changes described by the commit:
#  - add a helper average (col J, row 12)
#  - add a small summary block (rows 14-17) with labels + aggregate formulas
#    over "SW(S*)/SW(OPT)" (col N) and "SC(S*)/SC(OPT)" (col Z)
#  - summary values get a bold/size-12 font with vertical-center alignment

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: quick average of the |S*|/n column (J)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Row 14-17: summary labels (col A) + aggregate formulas (col B)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style the first summary cell: bold, size 12, vertically centered ...
$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108

# ... then replicate the same formatting onto the rest of the block via a
# format-only copy/paste so we don't spawn extra (unused) cell styles.
$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)

# Leave the selection where the author last worked (the helper-average cell)
$ws.Range("J12").Select() | Out-Null
